$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.560.17"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "2.284.08"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'303.93"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").Value = "'95.54"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("E7").Value = "  -2.77%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -2.99%  "
$ws.Range("D10").Value = "'34.71"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.41%  "
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").Value = "2.639.44"
$ws.Range("D16").Value = "2.292.10"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").Value = "'0.770"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("D18").Value = "42.471.25"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "'12.90"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").Value = "0.0₃0891"
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("E21").Value = "  -2.59%  "
$ws.Range("D22").Value = "'67.10"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.85%  "
$ws.Range("D23").Value = "'235.78"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.71%  "
$ws.Range("D24").Value = "'2.13"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -1.61%  "
$ws.Range("D27").Value = "'24.66"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("E28").Value = "  +17.17%  "
$ws.Range("D29").Value = "'166.40"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").Value = "'8.96"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("D31").Value = "'32.53"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.02%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "'17.81"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").Value = "'4.41"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -8.00%  "
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("E39").Value = "  -2.42%  "
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("E41").Value = "  -4.09%  "
$ws.Range("D42").Value = "1.989.24"
$ws.Range("E42").Value = "  -0.73%  "
$ws.Range("D43").Value = "'0.0276"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.76%  "
$ws.Range("D44").Value = "'18.40"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.61%  "
$ws.Range("D45").Value = "'10.16"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("D46").Value = "'2.00"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -7.79%  "
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("D48").Value = "'2.90"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.24%  "
$ws.Range("D49").Value = "'53.54"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "2.506.04"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("E51").Value = "  +0.50%  "
